$d = $word.ActiveDocument
$enDash = [char]0x2013
$brk = [char]11

# ---------------------------------------------------------------------------
# 1. "- Worker Management: Add/Edit, Toggle Availability"
#    -> "- " (plain) + "Worker Management: Add/Edit, Toggle Availability" (strike)
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Worker Management: Add/Edit, Toggle Availability", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Font.StrikeThrough = 1

# ---------------------------------------------------------------------------
# 2. "- Live Worker Map with Ping Status" -> "- " (plain) + text (strike)
#    the following line break becomes struck too, and merges with the
#    start of the "Dispatch Panel" line (its own <w:br/> disappears);
#    "- Dispatch Panel..." stays unstruck but is split by the _GoBack
#    bookmark after "- Dis".
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Live Worker Map with Ping Status", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rngWithBreak = $d.Range($rng.Start, $rng.End + 1)
$rngWithBreak.Font.StrikeThrough = 1

# ---------------------------------------------------------------------------
# 3. "- Route Replay for GPS History" -> "- " (plain) + text (strike)
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Route Replay for GPS History", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Font.StrikeThrough = 1

# ---------------------------------------------------------------------------
# 4. "- Category Manager - Add/Edit Main Category " -> "- " (plain) + text (strike)
# ---------------------------------------------------------------------------
$rng = $d.Content
$searchText = "Category Manager $enDash Add/Edit Main Category "
$rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Font.StrikeThrough = 1

# ---------------------------------------------------------------------------
# 4b. Merge the separate "Service" / " Manager - Add/Edit Main " / "Service"
#     runs (plus the lone <w:br/> run before them merging with the following
#     "- ") into "- " (plain) + "Service Manager - Add/Edit Main Service" (strike)
# ---------------------------------------------------------------------------
$rng = $d.Content
$serviceSearch = "Service Manager $enDash Add/Edit Main Service"
$rng.Find.Execute($serviceSearch, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$serviceStart = $rng.Start

# merge the lone <w:br/> run with the following "- " run by rewriting that
# span (break char + dash + space) through a placeholder first, since
# setting identical text is a no-op in this engine.
$prefixLen = 3
$prefixRng = $d.Range($serviceStart - $prefixLen, $serviceStart)
$prefixRng.Text = "ZZZ_PLACEHOLDER_ZZZ"
$prefixRng2 = $d.Range($serviceStart - $prefixLen, ($serviceStart - $prefixLen) + "ZZZ_PLACEHOLDER_ZZZ".Length)
$prefixRng2.Text = $brk + "- "

# re-find the "Service Manager..." text (positions may have shifted) and
# strike it - this also merges the split runs into a single run.
$rng2 = $d.Content
$rng2.Find.Execute($serviceSearch, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng2.Font.StrikeThrough = 1

# Move the _GoBack bookmark from the end of this line to between "- Dis" and
# "patch Panel..." on the earlier "Dispatch Panel" line.
$dispatchRng = $d.Content
$dispatchRng.Find.Execute("- Dispatch Panel: Drag and Drop or Auto-Assign", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bmPos = $dispatchRng.Start + ("- Dis").Length
$bmRng = $d.Range($bmPos, $bmPos)
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$d.Bookmarks.Add("_GoBack", $bmRng) | Out-Null

# ---------------------------------------------------------------------------
# 5. "Product Library - Add Products with Image, Price, and Assign to
#    Services" (bold) -> add strike (keep bold)
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Product Library $enDash Add Products with Image, Price, and Assign to Services", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Font.StrikeThrough = 1

# ---------------------------------------------------------------------------
# 6. "Coupon Manager - Create Flat/Percent Discounts, Usage Limits" (bold)
#    -> add strike (keep bold); the following " " run also gets strike
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Coupon Manager $enDash Create Flat/Percent Discounts, Usage Limits", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Font.StrikeThrough = 1
$spaceRng = $d.Range($rng.End, $rng.End + 1)
$spaceRng.Font.StrikeThrough = 1

# ---------------------------------------------------------------------------
# 7. "Ad Manager - Upload Banners, Target by Language/Nationality" (bold)
#    -> add strike (keep bold); following <w:br/> also struck; then
#    "- Pricing Settings: Base, Add-on, Delivery Fees" splits into "- "
#    (plain) + text (strike)
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Ad Manager $enDash Upload Banners, Target by Language/Nationality", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Font.StrikeThrough = 1
$brRng = $d.Range($rng.End, $rng.End + 1)
$brRng.Font.StrikeThrough = 1

$rng = $d.Content
$rng.Find.Execute("Pricing Settings: Base, Add-on, Delivery Fees", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Font.StrikeThrough = 1

# ---------------------------------------------------------------------------
# 8. "- Feedback Viewer - Customer Ratings & Comments" -> "- " (plain) +
#    text (strike)
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Feedback Viewer $enDash Customer Ratings & Comments", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Font.StrikeThrough = 1

# ---------------------------------------------------------------------------
# 9. "- Stripe (Global) - Apple Pay, Google Pay, Recurring Billing"
#    -> "- Stripe (Global) " (trims the rest)
# ---------------------------------------------------------------------------
$rng = $d.Content
$oldStripe = "- Stripe (Global) $enDash Apple Pay, Google Pay, Recurring Billing"
$newStripe = "- Stripe (Global) "
$rng.Find.Execute($oldStripe, $true, $false, $false, $false, $false, $true, 1, $false, $newStripe, 2) | Out-Null
